# Applies the Project Planner updates described by the commit:
# "Updated project planner." — refreshes actual-start/duration and
# percent-complete figures for several tasks, plus the currently
# highlighted period and the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# Currently highlighted period (H2)
$ws.Range("H2").Value = 46

# Row 5 - percent complete
$ws.Range("G5").Value = 0.3

# Row 9 - actual duration + percent complete
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 1

# Row 11 - actual duration + percent complete
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 1

# Row 12 - actual start/duration + percent complete
$ws.Range("E12").Value = 44
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1

# Row 13 - actual start + percent complete
$ws.Range("E13").Value = 44
$ws.Range("G13").Value = 0.15

# Row 15 - actual start
$ws.Range("E15").Value = 44

# Row 16 - actual start/duration + percent complete
$ws.Range("E16").Value = 44
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1

# Row 17 - actual start + percent complete
$ws.Range("E17").Value = 44
$ws.Range("G17").Value = 0.5

# Row 18 - actual start/duration + percent complete
$ws.Range("E18").Value = 44
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.9

# Row 19 - actual start
$ws.Range("E19").Value = 43

# Row 20 - actual start/duration + percent complete
$ws.Range("E20").Value = 43
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 1

# Row 21 - actual start/duration + percent complete
$ws.Range("E21").Value = 43
$ws.Range("F21").Value = 2
$ws.Range("G21").Value = 0.8

# Row 22 - actual duration + percent complete
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1

# Row 23 - actual start
$ws.Range("E23").Value = 43

# Row 24 - actual start + percent complete
$ws.Range("E24").Value = 43
$ws.Range("G24").Value = 0.7

# Row 25 - actual start + percent complete
$ws.Range("E25").Value = 43
$ws.Range("G25").Value = 0.5

# Row 26 - actual start + percent complete
$ws.Range("E26").Value = 43
$ws.Range("G26").Value = 0.8

# Row 27 - percent complete
$ws.Range("G27").Value = 0.9

# Row 28 - percent complete
$ws.Range("G28").Value = 0.1

# Force recalculation so the AVERAGE() formulas in G10/G15/G19/G23 refresh
$excel.CalculateFullRebuild()

# Restore the active-cell selection recorded in the saved view state
$ws.Range("G5").Select()
